# Update Trade_Confirmation_extracted.xlsx:
#  1. Clear J2 ("Company Name or Bank Name" header value) on "Extracted Fields" sheet.
#  2. Add a new "Field Metadata" worksheet right after "Extracted Fields" describing
#     how each extracted field was parsed (byte offsets into the source doc, etc.)

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- 1. Clear J2 on the existing sheet ---
$ws1.Range("J2").Value = ""

# --- 2. Add the new "Field Metadata" sheet ---
$meta = $wb.Worksheets.Add($null, $ws1)
$meta.Name = "Field Metadata"

# Build the header format (bold font, centered/top-aligned, thin box border)
# on A1 first, then clone it onto B1:O1 so only a single extra cell style
# entry is added to the stylesheet.
$a1 = $meta.Range("A1")
$a1.Font.Bold = $true
$a1.HorizontalAlignment = -4108
$a1.VerticalAlignment = -4160
$a1.Borders.LineStyle = 1

$a1.Copy()
$meta.Range("B1:O1").PasteSpecial(-4122)
$meta.Application.CutCopyMode = $false

# Header row values
$meta.Range("A1").Value = "Field"
$meta.Range("B1").Value = "RowIndex"
$meta.Range("C1").Value = "start_index_nbr"
$meta.Range("D1").Value = "end_index_nbr"
$meta.Range("E1").Value = "row_adder_cnt"
$meta.Range("F1").Value = "col_adder_cnt"
$meta.Range("G1").Value = "param_ref_delim_txt"
$meta.Range("H1").Value = "param_value_pos_cd"
$meta.Range("I1").Value = "unit_price_pct_ind"
$meta.Range("J1").Value = "param_nm_occur_ind"
$meta.Range("K1").Value = "date_format_cd"
$meta.Range("L1").Value = "decimal_separator_cd"
$meta.Range("M1").Value = "param_def_value_txt"
$meta.Range("N1").Value = "derivation_col"
$meta.Range("O1").Value = "operations_seq"

# Data rows: Field, RowIndex, start_index_nbr, end_index_nbr, row_adder_cnt,
#            col_adder_cnt, param_ref_delim_txt, param_value_pos_cd,
#            unit_price_pct_ind, param_nm_occur_ind, date_format_cd,
#            decimal_separator_cd, param_def_value_txt, derivation_col, operations_seq

$meta.Range("A2").Value = "Trade Date"
$meta.Range("B2").Value = 1
$meta.Range("C2").Value = 295
$meta.Range("D2").Value = 311
$meta.Range("E2").Value = 0
$meta.Range("F2").Value = 1
$meta.Range("G2").Value = ":"
$meta.Range("H2").Value = "R"
$meta.Range("I2").Value = $false
$meta.Range("J2").Value = 1
$meta.Range("K2").Value = "Month DD, YYYY"

$meta.Range("A3").Value = "Settlement Date"
$meta.Range("B3").Value = 1
$meta.Range("C3").Value = 322
$meta.Range("D3").Value = 342
$meta.Range("E3").Value = 0
$meta.Range("F3").Value = 1
$meta.Range("G3").Value = ":"
$meta.Range("H3").Value = "R"
$meta.Range("I3").Value = $false
$meta.Range("J3").Value = 1
$meta.Range("K3").Value = "Month DD, YYYY"

$meta.Range("A4").Value = "Transaction Type"
$meta.Range("B4").Value = 1
$meta.Range("N4").Value = "Buyer and Seller"
$meta.Range("O4").Value = "inferred from roles"

$meta.Range("A5").Value = "Net Amount or Consideration Amount or Nominal Amount"
$meta.Range("B5").Value = 1
$meta.Range("C5").Value = 494
$meta.Range("D5").Value = 512
$meta.Range("E5").Value = 0
$meta.Range("F5").Value = 1
$meta.Range("G5").Value = ":"
$meta.Range("H5").Value = "R"
$meta.Range("I5").Value = $false
$meta.Range("J5").Value = 1
$meta.Range("L5").Value = "."

$meta.Range("A6").Value = "Unit Price or Price"
$meta.Range("B6").Value = 1
$meta.Range("C6").Value = 523
$meta.Range("D6").Value = 529
$meta.Range("E6").Value = 0
$meta.Range("F6").Value = 1
$meta.Range("G6").Value = ":"
$meta.Range("H6").Value = "R"
$meta.Range("I6").Value = $true
$meta.Range("J6").Value = 1
$meta.Range("L6").Value = "."

$meta.Range("A7").Value = "Units or Shares"
$meta.Range("B7").Value = 1

$meta.Range("A8").Value = "ISIN or RIC"
$meta.Range("B8").Value = 1
$meta.Range("C8").Value = 473
$meta.Range("D8").Value = 483
$meta.Range("E8").Value = 0
$meta.Range("F8").Value = 1
$meta.Range("G8").Value = ":"
$meta.Range("H8").Value = "R"
$meta.Range("I8").Value = $false
$meta.Range("J8").Value = 1

$meta.Range("A9").Value = "Currency"
$meta.Range("B9").Value = 1
$meta.Range("C9").Value = 453
$meta.Range("D9").Value = 462
$meta.Range("E9").Value = 0
$meta.Range("F9").Value = 1
$meta.Range("G9").Value = ":"
$meta.Range("H9").Value = "R"
$meta.Range("I9").Value = $false
$meta.Range("J9").Value = 1

$meta.Range("A10").Value = "Company Name or Bank Name"
$meta.Range("B10").Value = 1
$meta.Range("N10").Value = "Buyer"

$meta.Range("A1").Select()
